# Apply crypto price/volume updates plus the Elrond/Algorand row swap
# (see commit: "Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.490.14'
$ws.Range('D3').Value = '1.844.79'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '262.76'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '0.5201'
$ws.Range('E7').Value = '  +1.14%  '
$ws.Range('D8').Value = '0.3222'
$ws.Range('E8').Value = '  -1.49%  '
$ws.Range('D9').Value = '0.06785'
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '18.65'
$ws.Range('E10').Value = '  -2.23%  '
$ws.Range('D11').Value = '0.7742'
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('D12').Value = '0.07771'
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('D13').Value = '1.853.20'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('D14').Value = '88.19'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.010'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('D17').Value = '13.92'
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007960'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('D20').Value = '26.533.39'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').Value = '2.087.36'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.610'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.19%  '
$ws.Range('D23').Value = '9.434'
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('D24').Value = '5.984'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('D25').Value = '143.11'
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('D26').Value = '2.163'
$ws.Range('E26').Value = '  -8.71%  '
$ws.Range('D27').Value = '1.684'
$ws.Range('E27').Value = '  +1.51%  '
$ws.Range('D28').Value = '16.98'
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('D29').Value = '111.48'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').Value = '4.157'
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('D31').Value = '0.08723'
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').Value = '4.097'
$ws.Range('E32').Value = '  -1.71%  '
$ws.Range('D33').Value = '0.04817'
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('D34').Value = '0.7191'
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('D35').Value = '1.126'
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.860'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('D37').Value = '3.089'
$ws.Range('E37').Value = '  -0.93%  '
$ws.Range('D38').Value = '0.01789'
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').Value = '2.199'
$ws.Range('E39').Value = '  -1.19%  '
$ws.Range('D40').Value = '0.4835'
$ws.Range('E40').Value = '  -1.88%  '
$ws.Range('D41').Value = '111.41'
$ws.Range('E41').Value = '  -1.85%  '
$ws.Range('D42').Value = '0.8899'
$ws.Range('E42').Value = '  -0.92%  '
$ws.Range('D43').Value = '6.032'
$ws.Range('E43').Value = '  -2.00%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').Value = '7.595'
$ws.Range('E45').Value = '  -2.78%  '
$ws.Range('D46').Value = '0.4176'
$ws.Range('E46').Value = '  -2.09%  '
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('D48').Value = '9.052'
$ws.Range('E48').Value = '  -1.36%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1230'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.08%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '34.87'
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('D51').Value = '0.8877'
$ws.Range('E51').Value = '  +3.88%  '
